# PVPlayTool DaS_Items.xlsx fix-up
#
# - "Lloyd's Talisman" texture filename had an apostrophe that kept the
#   texture from being found at runtime; the item's ImagePath cells are
#   corrected from tex_DaS_Lloyd'sTalisman.png -> tex_DaS_LloydsTalisman.png
# - The active selection / scroll position is moved down to the bottom of
#   the table (around row 43), matching where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the misspelled (apostrophe-containing) texture file name for every
# "Lloyd's Talisman" row so the texture loader can find the file.
$ws.Range("D43").Value = "tex_DaS_LloydsTalisman.png"
$ws.Range("D44").Value = "tex_DaS_LloydsTalisman.png"
$ws.Range("D45").Value = "tex_DaS_LloydsTalisman.png"

# Move the view / selection down near the bottom of the table and select
# D43, matching where editing left off.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 1
$ws.Range("D43").Select()
